# Updates Betfair back/lay odds values per the commit diff.
# Each assignment below corresponds to one changed cell in the source XLSX
# (rows 3-27 of Sheet1), grouped by row for readability.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("J3").Value = 1.09
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 1.75
$ws.Range("AN3").Value = 2.72

# Row 4
$ws.Range("K4").Value = 3.6
$ws.Range("L4").Value = 1.39
$ws.Range("O4").Value = 1.33
$ws.Range("T4").Value = 1.86

# Row 5
$ws.Range("F5").Value = 1.76
$ws.Range("H5").Value = 5.5
$ws.Range("J5").Value = 3.35
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 2.72
$ws.Range("O5").Value = 1.41
$ws.Range("T5").Value = 2.14
$ws.Range("AN5").Value = 23

# Row 6
$ws.Range("G6").Value = 2.76
$ws.Range("I6").Value = 2.94
$ws.Range("J6").Value = 3.8
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 2.78
$ws.Range("Q6").Value = 1.45
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 1.46
$ws.Range("W6").Value = 1.57

# Row 7
$ws.Range("I7").Value = 870
$ws.Range("K7").Value = 21
$ws.Range("T7").Value = 1.04
$ws.Range("U7").Value = 1.04

# Row 8
$ws.Range("Q8").Value = 1.6
$ws.Range("R8").Value = 1.63

# Row 9
$ws.Range("M9").Value = 1.06
$ws.Range("Q9").Value = 1.89

# Row 10
$ws.Range("N10").Value = 3.25

# Row 11
$ws.Range("N11").Value = 1.1
$ws.Range("T11").Value = 1.71
$ws.Range("U11").Value = 2.12
$ws.Range("V11").Value = 3.7

# Row 12
$ws.Range("G12").Value = 2.36
$ws.Range("Q12").Value = 1.71
$ws.Range("X12").Value = 18
$ws.Range("Y12").Value = 19.5
$ws.Range("AF12").Value = 16.5
$ws.Range("AN12").Value = 17.5
$ws.Range("AO12").Value = 25

# Row 13
$ws.Range("Q13").Value = 1.94
$ws.Range("R13").Value = 1.36
$ws.Range("X13").Value = 15
$ws.Range("AG13").Value = 9.800000000000001
$ws.Range("AN13").Value = 13

# Row 14
$ws.Range("G14").Value = 1.96
$ws.Range("K14").Value = 3.85
$ws.Range("U14").Value = 1.9

# Row 15
$ws.Range("H15").Value = 10
$ws.Range("K15").Value = 6.8
$ws.Range("W15").Value = 4
$ws.Range("AN15").Value = 4

# Row 16
$ws.Range("F16").Value = 3.4
$ws.Range("G16").Value = 3.7
$ws.Range("H16").Value = 2.14
$ws.Range("I16").Value = 2.26
$ws.Range("J16").Value = 3.6
$ws.Range("P16").Value = 2.02
$ws.Range("R16").Value = 1.39
$ws.Range("S16").Value = 3
$ws.Range("V16").Value = 1.79
$ws.Range("W16").Value = 1.37

# Row 17
$ws.Range("G17").Value = 3.2
$ws.Range("I17").Value = 2.58
$ws.Range("K17").Value = 4.5
$ws.Range("T17").Value = 1.52
$ws.Range("V17").Value = 1.64
$ws.Range("W17").Value = 1.46

# Row 18
$ws.Range("H18").Value = 1.98
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 5.2
$ws.Range("N18").Value = 7
$ws.Range("P18").Value = 3.4
$ws.Range("R18").Value = 1.98
$ws.Range("S18").Value = 1.75
$ws.Range("V18").Value = 1.83
$ws.Range("AC18").Value = 16
$ws.Range("AD18").Value = 15.5
$ws.Range("AN18").Value = 17.5
$ws.Range("AO18").Value = 8

# Row 19
$ws.Range("H19").Value = 2.26
$ws.Range("I19").Value = 2.62
$ws.Range("O19").Value = 1.5
$ws.Range("S19").Value = 4.7
$ws.Range("Z19").Value = 16.5
$ws.Range("AG19").Value = 18
$ws.Range("AH19").Value = 25

# Row 20
$ws.Range("N20").Value = 2.5
$ws.Range("P20").Value = 1.58
$ws.Range("Q20").Value = 1.99
$ws.Range("S20").Value = 2.82
$ws.Range("T20").Value = 1.04
$ws.Range("U20").Value = 1.04

# Row 21
$ws.Range("F21").Value = 1.52
$ws.Range("G21").Value = 1.53
$ws.Range("J21").Value = 4.3
$ws.Range("K21").Value = 4.4
$ws.Range("N21").Value = 3.55
$ws.Range("U21").Value = 1.71
$ws.Range("X21").Value = 12
$ws.Range("AK21").Value = 18
$ws.Range("AL21").Value = 48
$ws.Range("AN21").Value = 9.800000000000001

# Row 22
$ws.Range("K22").Value = 3.1

# Row 23
$ws.Range("I23").Value = 4.4
$ws.Range("N23").Value = 3.4
$ws.Range("O23").Value = 1.4
$ws.Range("Q23").Value = 2.22
$ws.Range("T23").Value = 1.96
$ws.Range("V23").Value = 1.29
$ws.Range("X23").Value = 11
$ws.Range("AG23").Value = 10.5
$ws.Range("AN23").Value = 17.5

# Row 24
$ws.Range("I24").Value = 1.58
$ws.Range("J24").Value = 4.5
$ws.Range("K24").Value = 4.6
$ws.Range("L24").Value = 1.31
$ws.Range("O24").Value = 1.29
$ws.Range("P24").Value = 2.06
$ws.Range("R24").Value = 1.4
$ws.Range("S24").Value = 3.2
$ws.Range("U24").Value = 1.93
$ws.Range("V24").Value = 2.72
$ws.Range("AD24").Value = 10

# Row 25
$ws.Range("F25").Value = 1.91
$ws.Range("G25").Value = 1.96
$ws.Range("H25").Value = 4.1
$ws.Range("I25").Value = 4.4
$ws.Range("J25").Value = 3.8
$ws.Range("P25").Value = 1.99
$ws.Range("Q25").Value = 1.91
$ws.Range("R25").Value = 1.37
$ws.Range("S25").Value = 3.25
$ws.Range("T25").Value = 1.79
$ws.Range("U25").Value = 2.14
$ws.Range("V25").Value = 1.29
$ws.Range("W25").Value = 2.04
$ws.Range("X25").Value = 16
$ws.Range("AA25").Value = 95
$ws.Range("AB25").Value = 9.4
$ws.Range("AC25").Value = 9.199999999999999
$ws.Range("AE25").Value = 55
$ws.Range("AF25").Value = 12.5
$ws.Range("AG25").Value = 10.5
$ws.Range("AH25").Value = 19
$ws.Range("AJ25").Value = 24
$ws.Range("AN25").Value = 14
$ws.Range("AO25").Value = 55

# Row 26
$ws.Range("G26").Value = 1.84
$ws.Range("H26").Value = 5.1
$ws.Range("S26").Value = 4.3
$ws.Range("T26").Value = 2.22
$ws.Range("AE26").Value = 150
$ws.Range("AM26").Value = 270

# Row 27
$ws.Range("F27").Value = 1.64
$ws.Range("T27").Value = 1.04
$ws.Range("U27").Value = 1.81
